$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 67

$ws.Range("C18").Value = 45
$ws.Range("D18").Value = 39

$ws.Range("D28").Value = 45

$ws.Range("D34").Value = 42

$ws.Range("C46").Value = 70
$ws.Range("D46").Value = 54

$ws.Range("C49").Value = 62

$ws.Range("C52").Value = 51
$ws.Range("D52").Value = 41

$ws.Range("C59").Value = 63
$ws.Range("D59").Value = 50

$ws.Range("C62").Value = 63
$ws.Range("D62").Value = 53

$ws.Range("C71").Value = 72
$ws.Range("D71").Value = 65

$ws.Range("D77").Value = 133

$ws.Range("C92").Value = 248
$ws.Range("D92").Value = 180

$ws.Range("C93").Value = 5472
$ws.Range("D93").Value = 4405
